$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.364.36'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '3.505.31'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '599.62'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = '175.87'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('E10').Value = '  -2.68%  '
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '4.115.20'
$ws.Range('D13').Value = '31.32'
$ws.Range('E13').Value = '  +9.85%  '
$ws.Range('D14').Value = '0.135'
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('D15').Value = '67.364.64'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '3.511.83'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').Value = '6.31'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').Value = '14.59'
$ws.Range('E19').Value = '  +2.92%  '
$ws.Range('D20').Value = '393.05'
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('D21').Value = '8.01'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '73.41'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').Value = '0.540'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '5.71'
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('D27').Value = '10.30'
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('D29').Value = '0.995'
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('E30').Value = '  -2.85%  '
$ws.Range('E31').Value = '  -2.91%  '
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('D33').Value = '23.71'
$ws.Range('E33').Value = '  -2.02%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  +1.45%  '
$ws.Range('D36').Value = '163.82'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +1.40%  '
$ws.Range('D38').Value = '0.879'
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').Value = '4.68'
$ws.Range('E40').Value = '  -2.08%  '
$ws.Range('D41').Value = '26.58'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').Value = '27.15'
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('D43').Value = '0.0732'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('D44').Value = '2.809.94'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '42.54'
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = '2.55'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('E47').Value = '  -3.94%  '
$ws.Range('D48').Value = '338.28'
$ws.Range('E48').Value = '  -1.21%  '
$ws.Range('E49').Value = '  -2.41%  '
$ws.Range('D50').Value = '33.62'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = '0.848'
$ws.Range('E51').Value = '  -0.69%  '
